# ---------------------------------------------------------------------------
# Adds the "ODI Bowling Extra" worksheet (mirroring the existing
# "ODI Batting Extra" sheet) and tidies up "ODI Batting Extra" by clearing
# the placeholder cells that never received a scraped value.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "ODI Bowling Extra" sheet right after "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

# Match the page geometry used by every other sheet in the workbook.
$bowlingExtra.PageSetup.LeftMargin   = 54
$bowlingExtra.PageSetup.RightMargin  = 54
$bowlingExtra.PageSetup.TopMargin    = 72
$bowlingExtra.PageSetup.BottomMargin = 72
$bowlingExtra.PageSetup.HeaderMargin = 36
$bowlingExtra.PageSetup.FooterMargin = 36

# Clone the bold/bordered/centered header formatting from the sheet next door
# so the new header row reuses the existing style instead of minting a new one.
$battingExtra.Range("A1").Copy() | Out-Null
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122) | Out-Null

# Header row
$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Every data column is textual (match codes, counts, and "xx.xx%" strings),
# matching the convention used throughout the rest of the workbook.
$bowlingExtraData = $bowlingExtra.Range("A2:C21")
$bowlingExtraData.NumberFormat = "@"

$rows = @(
    @("4465", "0", ""),
    @("4481", "0", "30.00%"),
    @("4537", "0", "30.00%"),
    @("4538", "", ""),
    @("4539", "0", ""),
    @("4550", "0", ""),
    @("4557", "0", ""),
    @("4559", "0", ""),
    @("4606", "", ""),
    @("4611", "0", ""),
    @("4616", "0", "20.00%"),
    @("4626", "0", "10.00%"),
    @("4628", "0", "40.00%"),
    @("4679", "1", ""),
    @("4682", "1", "10.00%"),
    @("4685", "", ""),
    @("4711", "0", ""),
    @("4713", "0", ""),
    @("4717", "0", "10.00%"),
    @("4726", "", "")
)

$r = 2
foreach ($row in $rows) {
    $bowlingExtra.Range("A$r").Value = $row[0]
    if ($row[1] -ne "") {
        $bowlingExtra.Range("B$r").Value = $row[1]
    }
    if ($row[2] -ne "") {
        $bowlingExtra.Range("C$r").Value = $row[2]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) "ODI Batting Extra": clear out the cells that only ever held an empty
#    placeholder (no scraped value) - leaves cells with real data untouched.
# ---------------------------------------------------------------------------
$battingExtra.Range("B2:E2").Value = ""
$battingExtra.Range("C4:E4").Value = ""
$battingExtra.Range("B6:E6").Value = ""
$battingExtra.Range("C7:E7").Value = ""
$battingExtra.Range("E8").Value = ""
$battingExtra.Range("C9:E9").Value = ""
$battingExtra.Range("B10:E10").Value = ""
$battingExtra.Range("E11").Value = ""
$battingExtra.Range("E12").Value = ""
$battingExtra.Range("B13:E13").Value = ""
$battingExtra.Range("B14:F14").Value = ""
$battingExtra.Range("B15:F15").Value = ""
$battingExtra.Range("B16:F16").Value = ""
$battingExtra.Range("B17:F17").Value = ""
$battingExtra.Range("B18:F18").Value = ""
$battingExtra.Range("B19:F19").Value = ""
$battingExtra.Range("B20:F20").Value = ""
$battingExtra.Range("B21:F21").Value = ""
